# Auto-generated Excel COM-interop edit script
# Applies the numeric profit/price recalculations described in the commit diff
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.
$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Cells.Item(19, 8).Value = 4983.4287
$ws.Cells.Item(19, 9).Value = 2500.4443
$ws.Cells.Item(19, 11).Value = 2500.4443
$ws.Cells.Item(19, 13).Value = -2325.4443
# Row 40
$ws.Cells.Item(40, 8).Value = 3200.4285
$ws.Cells.Item(40, 9).Value = 2950.25
$ws.Cells.Item(40, 10).Value = 3534
$ws.Cells.Item(40, 11).Value = 2950.25
$ws.Cells.Item(40, 12).Value = 3534
$ws.Cells.Item(40, 13).Value = -2775.25
$ws.Cells.Item(40, 14).Value = -3884
# Row 43
$ws.Cells.Item(43, 8).Value = 14690.909
$ws.Cells.Item(43, 10).Value = 12085.714
$ws.Cells.Item(43, 12).Value = 12085.714
$ws.Cells.Item(43, 14).Value = -12223.714
# Row 137
$ws.Cells.Item(137, 8).Value = 4713.3335
$ws.Cells.Item(137, 9).Value = 2094.85
$ws.Cells.Item(137, 11).Value = 6284.549999999999
$ws.Cells.Item(137, 13).Value = -3734.549999999999
# Row 138
$ws.Cells.Item(138, 8).Value = 3494.3977
$ws.Cells.Item(138, 9).Value = 2788.3333
$ws.Cells.Item(138, 10).Value = 3834.8215
$ws.Cells.Item(138, 11).Value = 8364.999899999999
$ws.Cells.Item(138, 12).Value = 11504.4645
$ws.Cells.Item(138, 13).Value = -3224.999899999999
$ws.Cells.Item(138, 14).Value = -21784.4645

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 21096.197
$ws.Cells.Item(32, 9).Value = 21078.383
$ws.Cells.Item(32, 10).Value = 21500
$ws.Cells.Item(32, 11).Value = 21078.383
$ws.Cells.Item(32, 12).Value = 21500
$ws.Cells.Item(32, 13).Value = -20791.383
$ws.Cells.Item(32, 14).Value = -22074
# Row 61
$ws.Cells.Item(61, 8).Value = 8234.120000000001
$ws.Cells.Item(61, 9).Value = 5999
$ws.Cells.Item(61, 11).Value = 5999
$ws.Cells.Item(61, 13).Value = -5787
# Row 63
$ws.Cells.Item(63, 8).Value = 9046.677
$ws.Cells.Item(63, 9).Value = 4948.8
$ws.Cells.Item(63, 10).Value = 9753.207
$ws.Cells.Item(63, 11).Value = 4948.8
$ws.Cells.Item(63, 12).Value = 9753.207
$ws.Cells.Item(63, 13).Value = -4262.8
$ws.Cells.Item(63, 14).Value = -11125.207
# Row 66
$ws.Cells.Item(66, 8).Value = 9046.677
$ws.Cells.Item(66, 9).Value = 4948.8
$ws.Cells.Item(66, 10).Value = 9753.207
$ws.Cells.Item(66, 11).Value = 24744
$ws.Cells.Item(66, 12).Value = 48766.035
$ws.Cells.Item(66, 13).Value = -21312
$ws.Cells.Item(66, 14).Value = -55630.035
# Row 74
$ws.Cells.Item(74, 8).Value = 348597
$ws.Cells.Item(74, 9).Value = 386127.47
$ws.Cells.Item(74, 11).Value = 386127.47
$ws.Cells.Item(74, 13).Value = -385253.47
# Row 77
$ws.Cells.Item(77, 8).Value = 348597
$ws.Cells.Item(77, 9).Value = 386127.47
$ws.Cells.Item(77, 11).Value = 1930637.35
$ws.Cells.Item(77, 13).Value = -1926269.35
# Row 122
$ws.Cells.Item(122, 8).Value = 45532.582
$ws.Cells.Item(122, 9).Value = 4121.3335
$ws.Cells.Item(122, 11).Value = 12364.0005
$ws.Cells.Item(122, 13).Value = -9914.000499999998
# Row 132
$ws.Cells.Item(132, 8).Value = 9736.103999999999
$ws.Cells.Item(132, 9).Value = 7473
$ws.Cells.Item(132, 10).Value = 15676.75
$ws.Cells.Item(132, 11).Value = 22419
$ws.Cells.Item(132, 12).Value = 47030.25
$ws.Cells.Item(132, 13).Value = -19889
$ws.Cells.Item(132, 14).Value = -52090.25
# Row 136
$ws.Cells.Item(136, 8).Value = 8234.120000000001
$ws.Cells.Item(136, 9).Value = 5999
$ws.Cells.Item(136, 11).Value = 17997
$ws.Cells.Item(136, 13).Value = -15447

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
# Row 97
$ws.Cells.Item(97, 8).Value = 5159.4
$ws.Cells.Item(97, 10).Value = 0
$ws.Cells.Item(97, 12).Value = 0
$ws.Cells.Item(97, 14).ClearContents()
# Row 107
$ws.Cells.Item(107, 8).Value = 1502.7333
$ws.Cells.Item(107, 9).Value = 1323.3846
$ws.Cells.Item(107, 11).Value = 1323.3846
$ws.Cells.Item(107, 13).Value = 596.6153999999999
# Row 134
$ws.Cells.Item(134, 8).Value = 5175.268
$ws.Cells.Item(134, 9).Value = 3930.9773
$ws.Cells.Item(134, 10).Value = 9737.666999999999
$ws.Cells.Item(134, 11).Value = 11792.9319
$ws.Cells.Item(134, 12).Value = 29213.001
$ws.Cells.Item(134, 13).Value = -9257.9319
$ws.Cells.Item(134, 14).Value = -34283.001

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
# Row 52
$ws.Cells.Item(52, 8).Value = 83499.5
$ws.Cells.Item(52, 10).Value = 96999
$ws.Cells.Item(52, 12).Value = 96999
$ws.Cells.Item(52, 14).Value = -97587
# Row 99
$ws.Cells.Item(99, 8).Value = 2351.7222
$ws.Cells.Item(99, 9).Value = 2483.1
$ws.Cells.Item(99, 10).Value = 2187.5
$ws.Cells.Item(99, 11).Value = 2483.1
$ws.Cells.Item(99, 12).Value = 2187.5
$ws.Cells.Item(99, 13).Value = -985.0999999999999
$ws.Cells.Item(99, 14).Value = -5183.5
# Row 106
$ws.Cells.Item(106, 8).Value = 73000
$ws.Cells.Item(106, 10).Value = 73000
$ws.Cells.Item(106, 12).Value = 73000
$ws.Cells.Item(106, 14).Value = -75524
# Row 126
$ws.Cells.Item(126, 8).Value = 2351.7222
$ws.Cells.Item(126, 9).Value = 2483.1
$ws.Cells.Item(126, 10).Value = 2187.5
$ws.Cells.Item(126, 11).Value = 7449.299999999999
$ws.Cells.Item(126, 12).Value = 6562.5
$ws.Cells.Item(126, 13).Value = -4979.299999999999
$ws.Cells.Item(126, 14).Value = -11502.5
# Row 132
$ws.Cells.Item(132, 8).Value = 17678.383
$ws.Cells.Item(132, 9).Value = 855.2
$ws.Cells.Item(132, 11).Value = 2565.6
$ws.Cells.Item(132, 13).Value = -35.60000000000036

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Cells.Item(2, 8).Value = 24348.2
$ws.Cells.Item(2, 10).Value = 40480.332
$ws.Cells.Item(2, 12).Value = 242881.992
$ws.Cells.Item(2, 14).Value = -243107.992
# Row 7
$ws.Cells.Item(7, 8).Value = 986071.3
$ws.Cells.Item(7, 9).Value = 1150349.9
$ws.Cells.Item(7, 11).Value = 3451049.7
$ws.Cells.Item(7, 13).Value = -3450937.7
# Row 34
$ws.Cells.Item(34, 8).Value = 2338
$ws.Cells.Item(34, 10).Value = 3300.125
$ws.Cells.Item(34, 12).Value = 9900.375
$ws.Cells.Item(34, 14).Value = -10068.375
# Row 131
$ws.Cells.Item(131, 8).Value = 10758522
$ws.Cells.Item(131, 10).Value = 8022.2383
$ws.Cells.Item(131, 12).Value = 24066.7149
$ws.Cells.Item(131, 14).Value = -34146.7149
# Row 141
$ws.Cells.Item(141, 8).Value = 4889.4287
$ws.Cells.Item(141, 9).Value = 4787.6665
$ws.Cells.Item(141, 10).Value = 5500
$ws.Cells.Item(141, 11).Value = 14362.9995
$ws.Cells.Item(141, 12).Value = 16500
$ws.Cells.Item(141, 13).Value = -9182.999500000002
$ws.Cells.Item(141, 14).Value = -26860

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
# Row 36
$ws.Cells.Item(36, 8).Value = 1665
$ws.Cells.Item(36, 9).Value = 1497.5
$ws.Cells.Item(36, 10).Value = 2000
$ws.Cells.Item(36, 11).Value = 1497.5
$ws.Cells.Item(36, 12).Value = 2000
$ws.Cells.Item(36, 13).Value = -1012.5
$ws.Cells.Item(36, 14).Value = -2970
# Row 97
$ws.Cells.Item(97, 8).Value = 2486.6667
$ws.Cells.Item(97, 9).Value = 2667.5
$ws.Cells.Item(97, 11).Value = 2667.5
$ws.Cells.Item(97, 13).Value = -2171.5

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 3779.0625
$ws.Cells.Item(22, 9).Value = 2182.7273
$ws.Cells.Item(22, 10).Value = 4615.2383
$ws.Cells.Item(22, 11).Value = 2182.7273
$ws.Cells.Item(22, 12).Value = 4615.2383
$ws.Cells.Item(22, 13).Value = -1887.7273
$ws.Cells.Item(22, 14).Value = -5205.2383
# Row 27
$ws.Cells.Item(27, 8).Value = 3779.0625
$ws.Cells.Item(27, 9).Value = 2182.7273
$ws.Cells.Item(27, 10).Value = 4615.2383
$ws.Cells.Item(27, 11).Value = 2182.7273
$ws.Cells.Item(27, 12).Value = 4615.2383
$ws.Cells.Item(27, 13).Value = -2075.7273
$ws.Cells.Item(27, 14).Value = -4829.2383
# Row 39
$ws.Cells.Item(39, 8).Value = 23166.334
$ws.Cells.Item(39, 9).Value = 4499
$ws.Cells.Item(39, 10).Value = 32500
$ws.Cells.Item(39, 11).Value = 4499
$ws.Cells.Item(39, 12).Value = 32500
$ws.Cells.Item(39, 13).Value = -4039
$ws.Cells.Item(39, 14).Value = -33420
# Row 55
$ws.Cells.Item(55, 8).Value = 839.53845
$ws.Cells.Item(55, 9).Value = 339.125
$ws.Cells.Item(55, 10).Value = 1640.2
$ws.Cells.Item(55, 11).Value = 339.125
$ws.Cells.Item(55, 12).Value = 1640.2
$ws.Cells.Item(55, 13).Value = -166.125
$ws.Cells.Item(55, 14).Value = -1986.2

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
# Row 19
$ws.Cells.Item(19, 8).Value = 0
$ws.Cells.Item(19, 10).Value = 0
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 14).ClearContents()

